$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows 6 and 7 - names
$ws.Range("A6").Value = "Nome 5"
$ws.Range("A7").Value = "Nome 6"

# New column C - "endereco" header plus values for rows 2-7
$ws.Range("C1").Value = "endereco"
$ws.Range("C2").Value = "rua x"
$ws.Range("C3").Value = "rua y"
$ws.Range("C4").Value = "rua z"
$ws.Range("C5").Value = "rua teste"
$ws.Range("C6").Value = "rua nova"
$ws.Range("C7").Value = "rua 123"

# New phone numbers for rows 6-7 (stored as numbers, not text)
$ws.Range("B6").Value = 4299999916
$ws.Range("B7").Value = 4299999917

# Update selection to C7
$ws.Range("C7").Select()
